# Updated Global Glider Cal and Ingest sheets
# - CC_scattering_angle (F2) changed to 140
# - CC_angular_resolution (F4) changed to 1.13
# - Selection left on F4 of Asset_Cal_Info sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Asset_Cal_Info")

$ws.Range("F2").Value = 140
$ws.Range("F4").Value = 1.13

$ws.Activate() | Out-Null
$ws.Range("F4").Select() | Out-Null
